$d = $word.ActiveDocument

# 1. Insert a new empty paragraph right after the paragraph ending
#    "... iff it is extensible." (and before the one with the page break).
#    Doing this via Find/Replace with a "^p" marker keeps the new
#    paragraph a clean <w:p/> with no stray run.
$taketypeParams = "Generated class should take type parameters R and T (extends source class) iff it is extensible."
$d.Content.Find.Execute(
    $taketypeParams,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $taketypeParams + "^p",
    2)

# 2. Collapse the five separate runs of the "Generated class should extend
#    com.mistraltech..." paragraph into a single run. Re-finding and
#    replacing the full paragraph text with itself makes Word rewrite the
#    match as one uniform run.
$oldExtend = "Generated class should extend com.mistraltech.smog.core.CompositePropertyMatcher if no superclass specified, passing R and T if extensible and this type and the source class if not."
$d.Content.Find.Execute($oldExtend, $true, $false, $false, $false, $false, $true, 1, $false, $oldExtend, 2)

# 3. Add a new sentence to the paragraph that follows it (the one holding
#    the _GoBack bookmark), inserting the text immediately before the
#    bookmark markers.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq $oldExtend) {
        $bookmarkPara = $p.Next()
        $r = $bookmarkPara.Range.Duplicate
        $r.Collapse(1)
        $r.InsertBefore("Generation works for properties with primitive and non-primitive types.")
        break
    }
}
